$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 457.0625
$ws.Range("I19").Value = 534.05884
$ws.Range("J19").Value = 369.8
$ws.Range("K19").Value = 534.05884
$ws.Range("L19").Value = 369.8
$ws.Range("M19").Value = -359.05884
$ws.Range("N19").Value = -719.8

$ws.Range("H76").Value = 6250
$ws.Range("I76").Value = 5000
$ws.Range("J76").Value = 7500
$ws.Range("K76").Value = 5000
$ws.Range("L76").Value = 7500
$ws.Range("M76").Value = -4685
$ws.Range("N76").Value = -8130

$ws.Range("H79").Value = 6250
$ws.Range("I79").Value = 5000
$ws.Range("J79").Value = 7500
$ws.Range("K79").Value = 5000
$ws.Range("L79").Value = 7500
$ws.Range("M79").Value = -3908
$ws.Range("N79").Value = -9684

$ws.Range("H101").Value = 573
$ws.Range("J101").Value = 392.8
$ws.Range("L101").Value = 1178.4
$ws.Range("N101").Value = -4422.4

$ws.Range("H113").Value = 1466.2222
$ws.Range("I113").Value = 1465
$ws.Range("K113").Value = 1465
$ws.Range("M113").Value = 1789

$ws.Range("H116").Value = 2490.125
$ws.Range("I116").Value = 2651.6667
$ws.Range("K116").Value = 2651.6667
$ws.Range("M116").Value = 790.3332999999998

$ws.Range("H138").Value = 2194.8
$ws.Range("J138").Value = 2439.8
$ws.Range("L138").Value = 7319.400000000001
$ws.Range("N138").Value = -17599.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4590.8076
$ws.Range("I32").Value = 4590.8076
$ws.Range("K32").Value = 4590.8076
$ws.Range("M32").Value = -4303.8076

$ws.Range("H46").Value = 4054.6
$ws.Range("I46").Value = 3387
$ws.Range("J46").Value = 4499.6665
$ws.Range("K46").Value = 3387
$ws.Range("L46").Value = 4499.6665
$ws.Range("M46").Value = -3068
$ws.Range("N46").Value = -5137.6665

$ws.Range("H74").Value = 6332.222
$ws.Range("I74").Value = 5874.125
$ws.Range("K74").Value = 5874.125
$ws.Range("M74").Value = -5000.125

$ws.Range("H77").Value = 6332.222
$ws.Range("I77").Value = 5874.125
$ws.Range("K77").Value = 29370.625
$ws.Range("M77").Value = -25002.625

$ws.Range("H122").Value = 1083.2727
$ws.Range("I122").Value = 1083.2727
$ws.Range("K122").Value = 3249.8181
$ws.Range("M122").Value = -799.8181

$ws.Range("H132").Value = 2722.7778
$ws.Range("I132").Value = 1592.6666
$ws.Range("K132").Value = 4777.9998
$ws.Range("M132").Value = -2247.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 94997
$ws.Range("J57").Value = 94997
$ws.Range("L57").Value = 94997
$ws.Range("N57").Value = -96437

$ws.Range("H80").Value = 203.33333
$ws.Range("I80").Value = 148.75
$ws.Range("J80").Value = 223.18182
$ws.Range("K80").Value = 148.75
$ws.Range("L80").Value = 223.18182
$ws.Range("M80").Value = 849.25
$ws.Range("N80").Value = -2219.18182

$ws.Range("H83").Value = 203.33333
$ws.Range("I83").Value = 148.75
$ws.Range("J83").Value = 223.18182
$ws.Range("K83").Value = 743.75
$ws.Range("L83").Value = 1115.9091
$ws.Range("M83").Value = 4248.25
$ws.Range("N83").Value = -11099.9091

$ws.Range("H105").Value = 1503.7693
$ws.Range("I105").Value = 1486.3636
$ws.Range("K105").Value = 1486.3636
$ws.Range("M105").Value = 260.6364000000001

$ws.Range("H107").Value = 3677.923
$ws.Range("I107").Value = 1414.125
$ws.Range("J107").Value = 7300
$ws.Range("K107").Value = 1414.125
$ws.Range("L107").Value = 7300
$ws.Range("M107").Value = 505.875
$ws.Range("N107").Value = -11140

$ws.Range("H134").Value = 4689.9287
$ws.Range("I134").Value = 3766.1
$ws.Range("K134").Value = 11298.3
$ws.Range("M134").Value = -8763.299999999999

$ws.Range("H136").Value = 94997
$ws.Range("J136").Value = 94997
$ws.Range("L136").Value = 94997
$ws.Range("N136").Value = -105197

$ws.Range("H140").Value = 106926
$ws.Range("J140").Value = 120780
$ws.Range("L140").Value = 120780
$ws.Range("N140").Value = -131140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1334.6666
$ws.Range("I3").Value = 1500.6666
$ws.Range("J3").Value = 1168.6666
$ws.Range("K3").Value = 1500.6666
$ws.Range("L3").Value = 1168.6666
$ws.Range("M3").Value = -1387.6666
$ws.Range("N3").Value = -1394.6666

$ws.Range("H7").Value = 74.70587999999999
$ws.Range("I7").Value = 36.615383
$ws.Range("J7").Value = 198.5
$ws.Range("K7").Value = 36.615383
$ws.Range("L7").Value = 198.5
$ws.Range("M7").Value = 76.38461699999999
$ws.Range("N7").Value = -424.5

$ws.Range("H16").Value = 1235.2632
$ws.Range("I16").Value = 1141.1666
$ws.Range("K16").Value = 1141.1666
$ws.Range("M16").Value = -854.1666

$ws.Range("H31").Value = 4652.921
$ws.Range("J31").Value = 9664.643
$ws.Range("L31").Value = 9664.643
$ws.Range("N31").Value = -10254.643

$ws.Range("H34").Value = 4652.921
$ws.Range("J34").Value = 9664.643
$ws.Range("L34").Value = 9664.643
$ws.Range("N34").Value = -10068.643

$ws.Range("H62").Value = 3133
$ws.Range("I62").Value = 3133
$ws.Range("K62").Value = 3133
$ws.Range("M62").Value = -2509

$ws.Range("H65").Value = 3133
$ws.Range("I65").Value = 3133
$ws.Range("K65").Value = 15665
$ws.Range("M65").Value = -12545

$ws.Range("H105").Value = 2079.3635
$ws.Range("I105").Value = 2079.3635
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2079.3635
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -332.3634999999999
$ws.Range("N105").ClearContents()

$ws.Range("H113").Value = 1235.2632
$ws.Range("I113").Value = 1141.1666
$ws.Range("K113").Value = 1141.1666
$ws.Range("M113").Value = 1028.8334

$ws.Range("H134").Value = 1491
$ws.Range("I134").Value = 1316.3889
$ws.Range("K134").Value = 3949.1667
$ws.Range("M134").Value = -1414.1667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2247.375
$ws.Range("I80").Value = 2244.75
$ws.Range("J80").Value = 2250
$ws.Range("K80").Value = 2244.75
$ws.Range("L80").Value = 2250
$ws.Range("M80").Value = -1246.75
$ws.Range("N80").Value = -4246

$ws.Range("H83").Value = 2247.375
$ws.Range("I83").Value = 2244.75
$ws.Range("J83").Value = 2250
$ws.Range("K83").Value = 11223.75
$ws.Range("L83").Value = 11250
$ws.Range("M83").Value = -6231.75
$ws.Range("N83").Value = -21234

$ws.Range("H102").Value = 1710.8
$ws.Range("J102").Value = 4998.3335
$ws.Range("L102").Value = 4998.3335
$ws.Range("N102").Value = -8242.333500000001

$ws.Range("H132").Value = 75640.69
$ws.Range("I132").Value = 97705.164
$ws.Range("J132").Value = 9447.25
$ws.Range("K132").Value = 293115.492
$ws.Range("L132").Value = 28341.75
$ws.Range("M132").Value = -290585.492
$ws.Range("N132").Value = -33401.75

$ws.Range("H137").Value = 139299.8
$ws.Range("J137").Value = 162124.75
$ws.Range("L137").Value = 162124.75
$ws.Range("N137").Value = -172324.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1190.4445
$ws.Range("I22").Value = 535.6667
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 535.6667
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = -240.6667
$ws.Range("N22").Value = -3090

$ws.Range("H27").Value = 1190.4445
$ws.Range("I27").Value = 535.6667
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 535.6667
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = -428.6667
$ws.Range("N27").Value = -2714

$ws.Range("H46").Value = 7417
$ws.Range("J46").Value = 8250.375
$ws.Range("L46").Value = 8250.375
$ws.Range("N46").Value = -8626.375

$ws.Range("H136").Value = 2660.353
$ws.Range("I136").Value = 2482.0667
$ws.Range("J136").Value = 3997.5
$ws.Range("K136").Value = 7446.2001
$ws.Range("L136").Value = 11992.5
$ws.Range("M136").Value = -4896.2001
$ws.Range("N136").Value = -17092.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H96").Value = 1448.4
$ws.Range("I96").Value = 1951.5
$ws.Range("J96").Value = 1113
$ws.Range("K96").Value = 1951.5
$ws.Range("L96").Value = 1113
$ws.Range("M96").Value = -578.5
$ws.Range("N96").Value = -3859

$ws.Range("H100").Value = 917.1667
$ws.Range("I100").Value = 1084.3334
$ws.Range("K100").Value = 2168.6668
$ws.Range("M100").Value = -1627.6668

$ws.Range("H126").Value = 2946.2666
$ws.Range("I126").Value = 1471.9474
$ws.Range("J126").Value = 5492.8184
$ws.Range("K126").Value = 4415.8422
$ws.Range("L126").Value = 16478.4552
$ws.Range("M126").Value = -1945.8422
$ws.Range("N126").Value = -21418.4552

$ws.Range("H136").Value = 3017.9312
$ws.Range("I136").Value = 1947.15
$ws.Range("K136").Value = 5841.450000000001
$ws.Range("M136").Value = -3291.450000000001
